# Generate Report for Handback
# Update handback-status.xlsx timestamps and status to reflect a newer
# handback report generation run.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for the two updated files
$wsOverview.Range("G2").Value = "2016-08-12 18:18:41"
$wsOverview.Range("G3").Value = "2016-08-12 18:18:41"

# zh-cn sheet: Priority changed from human translation (ht) to machine
# translation (mt), and the handoff / handback timestamps move forward.
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E3").Value = "mt"
$wsZhCn.Range("H2").Value = "2016-08-12 18:18:34"
$wsZhCn.Range("H3").Value = "2016-08-12 18:18:34"
$wsZhCn.Range("K2").Value = "2016-08-12 18:19:08"
$wsZhCn.Range("K3").Value = "2016-08-12 18:19:08"

# de-de sheet: Priority changed from human translation (ht) to machine
# translation (mt), handoff datetime matches the new Overview timestamp,
# and the handback datetime moves forward.
$wsDeDe.Range("E2").Value = "mt"
$wsDeDe.Range("E3").Value = "mt"
$wsDeDe.Range("H2").Value = "2016-08-12 18:18:41"
$wsDeDe.Range("H3").Value = "2016-08-12 18:18:41"
$wsDeDe.Range("K2").Value = "2016-08-12 18:19:17"
$wsDeDe.Range("K3").Value = "2016-08-12 18:19:17"
